# Update the carjacking-by-month-yoy-latest workbook with data through 2022-11-09.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and its tab to reflect the new "through" date.
$ws.Name = "Through 2022-11-09"

# Update the November row label.
$ws.Range("A12").Value = "November (through 11-09)"

# Updated November counts (row 12), by year column (2015..2022 => B..I).
$ws.Range("B12").Value = 12
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 32
$ws.Range("E12").Value = 23
$ws.Range("F12").Value = 13
$ws.Range("G12").Value = 57
$ws.Range("H12").Value = 67
$ws.Range("I12").Value = 25

# Updated Total counts (row 13), by year column (2015..2022 => B..I).
$ws.Range("B13").Value = 270
$ws.Range("C13").Value = 508
$ws.Range("D13").Value = 742
$ws.Range("E13").Value = 638
$ws.Range("F13").Value = 495
$ws.Range("G13").Value = 1114
$ws.Range("H13").Value = 1508
$ws.Range("I13").Value = 1424
